# Mailadressen.xlsx - update the address table:
#   - insert a new "Aaron Anneliese" row right after the header
#   - replace the old Mustermann/Mia/Roeckelein/Singer rows with the new
#     Caloni/Domini/Emil/Fan/Gans rows
#   - table grows from 6 rows (A1:C6) to 8 rows (A1:C8)
#
# Final table (Nachname | Vorname | Email):
#   1: Nachname  | Vorname  | Email
#   2: Aaron     | Anneliese| annelise@gmx.de
#   3: Blobfisch | Barbara  | Blobfisch@fischteich.de
#   4: Caloni    | Clara    | Caloni@gmx.de
#   5: Domini    | Dorothea | Domini@gmx.de
#   6: Emil      | Eva      | Domini@gmx.de
#   7: Fan       | Fabian   | Fabian@gmx.de
#   8: Gans      | Gustav   | Gans@gmx.de

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Nachname",  "Vorname",   "Email"),
    @("Aaron",     "Anneliese", "annelise@gmx.de"),
    @("Blobfisch", "Barbara",   "Blobfisch@fischteich.de"),
    @("Caloni",    "Clara",     "Caloni@gmx.de"),
    @("Domini",    "Dorothea",  "Domini@gmx.de"),
    @("Emil",      "Eva",       "Domini@gmx.de"),
    @("Fan",       "Fabian",    "Fabian@gmx.de"),
    @("Gans",      "Gustav",    "Gans@gmx.de")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 1
    $ws.Cells.Item($rowNum, 1).Value = $rows[$r][0]
    $ws.Cells.Item($rowNum, 2).Value = $rows[$r][1]
    $ws.Cells.Item($rowNum, 3).Value = $rows[$r][2]
}

$wb.Save()
